$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts D..K to E..L)
$ws.Columns("D").Insert()

# Copy formatting from column E (old D, now shifted) into new column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43100
$ws.Range("F7").Value = 42735
$ws.Range("G7").Value = 42369
$ws.Range("H7").Value = 42004
$ws.Range("I7").Value = 41639
$ws.Range("J7").Value = 41274
$ws.Range("K7").Value = 40908
$ws.Range("D8").Value = 10180000
$ws.Range("E8").Value = 10008900
$ws.Range("F8").Value = 9337600
$ws.Range("G8").Value = 8449000
$ws.Range("H8").Value = 7579800
$ws.Range("I8").Value = 6546800
$ws.Range("J8").Value = 5958900
$ws.Range("K8").Value = 6022200
$ws.Range("D9").Value = 3779800
$ws.Range("E9").Value = 4679600
$ws.Range("F9").Value = 4489800
$ws.Range("G9").Value = 4122900
$ws.Range("H9").Value = 3578900
$ws.Range("I9").Value = 2970100
$ws.Range("J9").Value = 2698500
$ws.Range("K9").Value = 2572800
$ws.Range("D10").Value = 6400200
$ws.Range("E10").Value = 5329300
$ws.Range("F10").Value = 4847800
$ws.Range("G10").Value = 4326100
$ws.Range("H10").Value = 4000900
$ws.Range("I10").Value = 3576700
$ws.Range("J10").Value = 3260400
$ws.Range("K10").Value = 3449400
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("D14").Value = 596100
$ws.Range("E14").Value = 935500
$ws.Range("F14").Value = 1114100
$ws.Range("G14").Value = 2787100
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = "NA"
$ws.Range("J14").Value = "NA"
$ws.Range("K14").Value = "NA"
$ws.Range("D15").Value = 1748900
$ws.Range("E15").Value = 1529400
$ws.Range("F15").Value = 1443000
$ws.Range("G15").Value = 1313500
$ws.Range("H15").Value = 1198900
$ws.Range("I15").Value = 1105200
$ws.Range("J15").Value = 1023300
$ws.Range("K15").Value = 1020200
$ws.Range("D17").Value = 6548300
$ws.Range("E17").Value = 7567900
$ws.Range("F17").Value = 7460000
$ws.Range("G17").Value = 8608300
$ws.Range("H17").Value = 5129800
$ws.Range("I17").Value = 4406500
$ws.Range("J17").Value = 4044800
$ws.Range("K17").Value = 3908000
$ws.Range("D18").Value = 3631700
$ws.Range("E18").Value = 2441000
$ws.Range("F18").Value = 1877600
$ws.Range("G18").Value = -159300
$ws.Range("H18").Value = 2449900
$ws.Range("I18").Value = 2140300
$ws.Range("J18").Value = 1914100
$ws.Range("K18").Value = 2114200
$ws.Range("D20").Value = 992800
$ws.Range("E20").Value = 1558400
$ws.Range("F20").Value = 146600
$ws.Range("G20").Value = 442100
$ws.Range("H20").Value = 543300
$ws.Range("I20").Value = 485200
$ws.Range("J20").Value = -307400
$ws.Range("K20").Value = 399500
$ws.Range("D21").Value = 6377100
$ws.Range("E21").Value = 5531900
$ws.Range("F21").Value = 3470300
$ws.Range("G21").Value = 1599100
$ws.Range("H21").Value = 4194600
$ws.Range("I21").Value = 3733000
$ws.Range("J21").Value = 2632200
$ws.Range("K21").Value = 3538700
$ws.Range("D22").Value = 1685600
$ws.Range("E22").Value = 1539000
$ws.Range("F22").Value = 1401300
$ws.Range("G22").Value = 1105900
$ws.Range("H22").Value = 891600
$ws.Range("I22").Value = 748700
$ws.Range("J22").Value = 164500
$ws.Range("K22").Value = 758200
$ws.Range("D23").Value = 2938900
$ws.Range("E23").Value = 2460400
$ws.Range("F23").Value = 622900
$ws.Range("G23").Value = -823100
$ws.Range("H23").Value = 2101600
$ws.Range("I23").Value = 1876900
$ws.Range("J23").Value = 1442300
$ws.Range("K23").Value = 1755400
$ws.Range("D24").Value = 445800
$ws.Range("E24").Value = 532100
$ws.Range("F24").Value = 262000
$ws.Range("G24").Value = 25300
$ws.Range("H24").Value = 618400
$ws.Range("I24").Value = 454700
$ws.Range("J24").Value = 346800
$ws.Range("K24").Value = 441700
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("D26").Value = 2493100
$ws.Range("E26").Value = 1928200
$ws.Range("F26").Value = 360900
$ws.Range("G26").Value = -848400
$ws.Range("H26").Value = 1483200
$ws.Range("I26").Value = 1422200
$ws.Range("J26").Value = 1095500
$ws.Range("K26").Value = 1313700
$ws.Range("D27").Value = 2509500
$ws.Range("E27").Value = 1632100
$ws.Range("F27").Value = 92300
$ws.Range("G27").Value = -922800
$ws.Range("H27").Value = 1297200
$ws.Range("I27").Value = 1274100
$ws.Range("J27").Value = 966700
$ws.Range("K27").Value = 1172300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("D29").Value = 124300
$ws.Range("E29").Value = 598300
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = "NA"
$ws.Range("I29").Value = "NA"
$ws.Range("J29").Value = "NA"
$ws.Range("K29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("D32").Value = -992800
$ws.Range("E32").Value = -1558400
$ws.Range("F32").Value = -146600
$ws.Range("G32").Value = -442100
$ws.Range("H32").Value = -543300
$ws.Range("I32").Value = -485200
$ws.Range("J32").Value = 307400
$ws.Range("K32").Value = -399500
$ws.Range("D33").Value = 2633800
$ws.Range("E33").Value = 2230400
$ws.Range("F33").Value = 92300
$ws.Range("G33").Value = -922800
$ws.Range("H33").Value = 1297200
$ws.Range("I33").Value = 1274100
$ws.Range("J33").Value = 966700
$ws.Range("K33").Value = 1172300
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("D35").Value = 2633800
$ws.Range("E35").Value = 2230400
$ws.Range("F35").Value = 92300
$ws.Range("G35").Value = -922800
$ws.Range("H35").Value = 1297200
$ws.Range("I35").Value = 1274100
$ws.Range("J35").Value = 966700
$ws.Range("K35").Value = 1172300
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43100
$ws.Range("F38").Value = 42735
$ws.Range("G38").Value = 42369
$ws.Range("H38").Value = 42004
$ws.Range("I38").Value = 41639
$ws.Range("J38").Value = 41274
$ws.Range("K38").Value = 40908
$ws.Range("D41").Value = 331900
$ws.Range("E41").Value = 810400
$ws.Range("F41").Value = 756100
$ws.Range("G41").Value = 632600
$ws.Range("H41").Value = 363900
$ws.Range("I41").Value = 689900
$ws.Range("J41").Value = 410100
$ws.Range("K41").Value = 502400
$ws.Range("D42").Value = 533600
$ws.Range("E42").Value = "NA"
$ws.Range("F42").Value = "NA"
$ws.Range("G42").Value = "NA"
$ws.Range("H42").Value = "NA"
$ws.Range("I42").Value = "NA"
$ws.Range("J42").Value = "NA"
$ws.Range("K42").Value = 0
$ws.Range("D43").Value = 2004900
$ws.Range("E43").Value = 1876900
$ws.Range("F43").Value = 1544200
$ws.Range("G43").Value = 1032200
$ws.Range("H43").Value = 977100
$ws.Range("I43").Value = 835000
$ws.Range("J43").Value = 782900
$ws.Range("K43").Value = 840400
$ws.Range("D44").Value = 320800
$ws.Range("E44").Value = 281300
$ws.Range("F44").Value = 273900
$ws.Range("G44").Value = 240400
$ws.Range("H44").Value = 217300
$ws.Range("I44").Value = 186800
$ws.Range("J44").Value = 166700
$ws.Range("K44").Value = 190500
$ws.Range("D45").Value = 630300
$ws.Range("E45").Value = 514200
$ws.Range("F45").Value = 3442000
$ws.Range("G45").Value = 1524100
$ws.Range("H45").Value = 1834500
$ws.Range("I45").Value = 630300
$ws.Range("J45").Value = 742000
$ws.Range("K45").Value = 855800
$ws.Range("D46").Value = 3821500
$ws.Range("E46").Value = 3482900
$ws.Range("F46").Value = 6016200
$ws.Range("G46").Value = 2915800
$ws.Range("H46").Value = 2316700
$ws.Range("I46").Value = 2342000
$ws.Range("J46").Value = 2101600
$ws.Range("K46").Value = 2389200
$ws.Range("D47").Value = 6309400
$ws.Range("E47").Value = 5472900
$ws.Range("F47").Value = 4957900
$ws.Range("G47").Value = 4836600
$ws.Range("H47").Value = 4344700
$ws.Range("I47").Value = 4505400
$ws.Range("J47").Value = 4278500
$ws.Range("K47").Value = 4236000
$ws.Range("D48").Value = 49492200
$ws.Range("E48").Value = 42626100
$ws.Range("F48").Value = 40540800
$ws.Range("G48").Value = 66706500
$ws.Range("H48").Value = 62177300
$ws.Range("I48").Value = 27986800
$ws.Range("J48").Value = 25089600
$ws.Range("K48").Value = 93656500
$ws.Range("D49").Value = 11333600
$ws.Range("E49").Value = 10180800
$ws.Range("F49").Value = 11946100
$ws.Range("G49").Value = 7278400
$ws.Range("H49").Value = 4364000
$ws.Range("I49").Value = 3416700
$ws.Range("J49").Value = 2853300
$ws.Range("K49").Value = 4594800
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("D52").Value = 2660600
$ws.Range("E52").Value = 2314500
$ws.Range("F52").Value = 2067400
$ws.Range("G52").Value = 1791300
$ws.Range("H52").Value = 1643200
$ws.Range("I52").Value = 1860500
$ws.Range("J52").Value = 1646900
$ws.Range("K52").Value = 1738500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("D54").Value = 73617300
$ws.Range("E54").Value = 64077200
$ws.Range("F54").Value = 65528400
$ws.Range("G54").Value = 47925600
$ws.Range("H54").Value = 43554900
$ws.Range("I54").Value = 40111400
$ws.Range("J54").Value = 35969900
$ws.Range("K54").Value = 36366500
$ws.Range("D57").Value = 2898700
$ws.Range("E57").Value = 2731300
$ws.Range("F57").Value = 2054000
$ws.Range("G57").Value = 2721600
$ws.Range("H57").Value = 1208600
$ws.Range("I57").Value = 644500
$ws.Range("J57").Value = 686900
$ws.Range("K57").Value = 534700
$ws.Range("D58").Value = 4632000
$ws.Range("E58").Value = 3444900
$ws.Range("F58").Value = 1943900
$ws.Range("G58").Value = 2802000
$ws.Range("H58").Value = 3173300
$ws.Range("I58").Value = 2095000
$ws.Range("J58").Value = 2358400
$ws.Range("K58").Value = 2174900
$ws.Range("D59").Value = 2103900
$ws.Range("E59").Value = 1174400
$ws.Range("F59").Value = 1717600
$ws.Range("G59").Value = 1514500
$ws.Range("H59").Value = 1262200
$ws.Range("I59").Value = 1248000
$ws.Range("J59").Value = 1331400
$ws.Range("K59").Value = 1558000
$ws.Range("D60").Value = 9634500
$ws.Range("E60").Value = 7350600
$ws.Range("F60").Value = 5715500
$ws.Range("G60").Value = 5478900
$ws.Range("H60").Value = 5641100
$ws.Range("I60").Value = 3987500
$ws.Range("J60").Value = 4376700
$ws.Range("K60").Value = 4242200
$ws.Range("D61").Value = 32757900
$ws.Range("E61").Value = 28936400
$ws.Range("F61").Value = 31437700
$ws.Range("G61").Value = 23307200
$ws.Range("H61").Value = 17950300
$ws.Range("I61").Value = 17083300
$ws.Range("J61").Value = 14149700
$ws.Range("K61").Value = 14396600
$ws.Range("D62").Value = 8159500
$ws.Range("E62").Value = 7777700
$ws.Range("F62").Value = 9038400
$ws.Range("G62").Value = 6566200
$ws.Range("H62").Value = 4593300
$ws.Range("I62").Value = 4055200
$ws.Range("J62").Value = 3797700
$ws.Range("K62").Value = 4414200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("D66").Value = 51783600
$ws.Range("E66").Value = 45443000
$ws.Range("F66").Value = 47476100
$ws.Range("G66").Value = 35692300
$ws.Range("H66").Value = 29362800
$ws.Range("I66").Value = 26324900
$ws.Range("J66").Value = 23384600
$ws.Range("K66").Value = 23464800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("D70").Value = 2962000
$ws.Range("E70").Value = 2962000
$ws.Range("F70").Value = 2962000
$ws.Range("G70").Value = 1859800
$ws.Range("H70").Value = 1678200
$ws.Range("I70").Value = 1349300
$ws.Range("J70").Value = 910900
$ws.Range("K70").Value = 1880600
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("D72").Value = 2063700
$ws.Range("E72").Value = 1207900
$ws.Range("F72").Value = 846900
$ws.Range("G72").Value = 2060700
$ws.Range("H72").Value = 4076800
$ws.Range("I72").Value = 3792500
$ws.Range("J72").Value = 3488100
$ws.Range("K72").Value = 3555400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("D76").Value = 18871600
$ws.Range("E76").Value = 15672300
$ws.Range("F76").Value = 15090300
$ws.Range("G76").Value = 10373500
$ws.Range("H76").Value = 12513900
$ws.Range("I76").Value = 12437200
$ws.Range("J76").Value = 11674400
$ws.Range("K76").Value = 11021100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43100
$ws.Range("F80").Value = 42735
$ws.Range("G80").Value = 42369
$ws.Range("H80").Value = 42004
$ws.Range("I80").Value = 41639
$ws.Range("J80").Value = 41274
$ws.Range("K80").Value = 40908
$ws.Range("D81").Value = 2633800
$ws.Range("E81").Value = 2230400
$ws.Range("F81").Value = 92300
$ws.Range("G81").Value = -922800
$ws.Range("H81").Value = 1297200
$ws.Range("I81").Value = 1274100
$ws.Range("J81").Value = 966700
$ws.Range("K81").Value = 1172300
$ws.Range("D83").Value = 1748900
$ws.Range("E83").Value = 1529400
$ws.Range("F83").Value = 1443000
$ws.Range("G83").Value = 1313500
$ws.Range("H83").Value = 1198900
$ws.Range("I83").Value = 1105200
$ws.Range("J83").Value = 1023300
$ws.Range("K83").Value = 1020200
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("D89").Value = 4878300
$ws.Range("E89").Value = 3892200
$ws.Range("F89").Value = 3772400
$ws.Range("G89").Value = 3062400
$ws.Range("H89").Value = 3035600
$ws.Range("I89").Value = 2734200
$ws.Range("J89").Value = 2657600
$ws.Range("K89").Value = 2831700
$ws.Range("D91").Value = -7378100
$ws.Range("E91").Value = -5603200
$ws.Range("F91").Value = -3945800
$ws.Range("G91").Value = -3296100
$ws.Range("H91").Value = -3242500
$ws.Range("I91").Value = -3319900
$ws.Range("J91").Value = -1931200
$ws.Range("K91").Value = -1930600
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("D94").Value = -7456200
$ws.Range("E94").Value = -2752800
$ws.Range("F94").Value = -13978500
$ws.Range("G94").Value = -3430800
$ws.Range("H94").Value = -3084000
$ws.Range("I94").Value = -3810400
$ws.Range("J94").Value = -2423100
$ws.Range("K94").Value = -2346200
$ws.Range("D96").Value = -1286700
$ws.Range("E96").Value = -1111800
$ws.Range("F96").Value = -1143100
$ws.Range("G96").Value = -1144600
$ws.Range("H96").Value = -1070900
$ws.Range("I96").Value = -1009100
$ws.Range("J96").Value = -953300
$ws.Range("K96").Value = -1561000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("D100").Value = 2045100
$ws.Range("E100").Value = -1056000
$ws.Range("F100").Value = 10424100
$ws.Range("G100").Value = 553700
$ws.Range("H100").Value = -277600
$ws.Range("I100").Value = 1335100
$ws.Range("J100").Value = -299900
$ws.Range("K100").Value = -493200
$ws.Range("D101").Value = 54300
$ws.Range("E101").Value = -29000
$ws.Range("F101").Value = -94500
$ws.Range("G101").Value = 83400
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 20800
$ws.Range("J101").Value = -11200
$ws.Range("K101").Value = 3100
$ws.Range("D102").Value = -478500
$ws.Range("E102").Value = 54300
$ws.Range("F102").Value = 123500
$ws.Range("G102").Value = 268700
$ws.Range("H102").Value = -326000
$ws.Range("I102").Value = 279800
$ws.Range("J102").Value = -76700
$ws.Range("K102").Value = -4600
